# Automatische test-sync: 2025-08-03 18:15:50
# Append a new log row (row 30) to the "Logs" sheet and bump the
# "Planning / Afspraak" count on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A30").Value = "Wil je dit oppakken?"
$logs.Range("B30").Value = "mailmind.test@zohomail.eu"
$logs.Range("C30").Value = "Testmail #2: Wil je dit oppakken?"
$logs.Range("D30").Value = "Planning / Afspraak"
$logs.Range("E30").Value = "Bedankt, we hebben dit doorgestuurd naar planning@bedrijf.nl."
$logs.Range("F30").Value = "2025-08-03 18:14:51"
$logs.Range("G30").Value = "Ja"
$logs.Range("H30").Value = "Ja"
$logs.Range("I30").Value = "Nee"
$logs.Range("J30").Value = "Nee"


# Conditional formatting ranges were previously bound to the last data
# row (29); extend them to cover the newly appended row 30 as well.
$logs.Range("D2:D29").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D30"))
$logs.Range("G2:G29").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G30"))
$logs.Range("H2:H29").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H30"))
$logs.Range("I2:I29").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I30"))
$logs.Range("J2:J29").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J30"))

$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B3").Value = 8
